$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert new "recipe_type" column before the old "description" column (F) ---
# This shifts the old F ("description") data/format to G and keeps header styling.
$ws.Range("F1").EntireColumn.Insert()
$ws.Range("F1").Value = "recipe_type"
$ws.Range("G1").Value = "description"

# Header row (row 1) gets taller once the new column is in place
$ws.Rows("1:1").RowHeight = 30

# --- Row 2: Base Yaourt Brasse Vache Nature -> PROCESS recipe ---
$ws.Range("F2").Value = "PROCESS"
$ws.Range("A2").Formula = '="REC_"&LEFT(F2,4)&"_"&C2'

# --- Row 3: Base Yaourt Brasse Vache Sucre -> PROCESS recipe ---
$ws.Range("F3").Value = "PROCESS"
$ws.Range("A3").Formula = '="REC_"&LEFT(F3,4)&"_"&C3'

# --- Row 4 (new row): CONDITIONNEMENT recipe ---
$ws.Range("B4").Value2 = $ws.Range("B2").Value2
$ws.Range("C4").Value = 1025700
$ws.Range("D4").Value = 1000
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = "CONDITIONNEMENT"
$ws.Range("G4").Value = "YAOURT BRASSE AU LAIT DE VACHE NATURE X 25 CANTINE BIOCHAMPS"
$ws.Range("A4").Formula = '="REC_"&LEFT(F4,4)&"_"&C4'

# Carry the "article_output_code" cell formatting (wrap + vertical-center) down
# column C through row 9, matching the original C2 look, without inventing new styles.
$ws.Range("C2").Copy()
$ws.Range("C4:C9").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Column widths ---
$ws.Range("A1").ColumnWidth = 41.583333333333336
$ws.Range("F1").ColumnWidth = 18.083333333333336
$ws.Range("G1").ColumnWidth = 63

$null = $ws.Range("A1:G4").Select()
